$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = $origStyle
}

$ws.Range("D2").Value = "60.727.75"
$ws.Range("E2").Value = "  -1.98%  "
$ws.Range("D3").Value = "2.404.20"
$ws.Range("E3").Value = "  -1.86%  "
$ws.Range("E4").Value = "  -0.11%  "
Set-TextValue "D5" "563.92"
$ws.Range("E5").Value = "  -2.65%  "
Set-TextValue "D6" "136.58"
$ws.Range("E6").Value = "  -3.75%  "
$ws.Range("E7").Value = "  +0.29%  "
Set-TextValue "D8" "0.534"
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("D9").Value = "2.386.96"
$ws.Range("E9").Value = "  -2.35%  "
Set-TextValue "D10" "0.105"
$ws.Range("E10").Value = "  -3.60%  "
$ws.Range("E11").Value = "  -0.61%  "
Set-TextValue "D12" "5.01"
$ws.Range("E12").Value = "  -3.19%  "
Set-TextValue "D13" "0.334"
$ws.Range("E13").Value = "  -1.67%  "
Set-TextValue "D14" "25.57"
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("D15").Value = "2.843.34"
$ws.Range("E15").Value = "  -1.53%  "
Set-TextValue "D16" "0.0000166"
$ws.Range("E16").Value = "  -3.72%  "
$ws.Range("D17").Value = "60.861.96"
$ws.Range("E17").Value = "  -1.67%  "
$ws.Range("D18").Value = "2.397.18"
$ws.Range("E18").Value = "  -1.56%  "
Set-TextValue "D19" "7.90"
$ws.Range("E19").Value = "  +9.64%  "
Set-TextValue "D20" "10.43"
$ws.Range("E20").Value = "  -1.79%  "
Set-TextValue "D21" "321.86"
$ws.Range("E21").Value = "  -0.90%  "
Set-TextValue "D22" "4.01"
$ws.Range("E22").Value = "  -1.37%  "
Set-TextValue "D23" "6.12"
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("E24").Value = "  +0.06%  "
Set-TextValue "D25" "1.78"
$ws.Range("E25").Value = "  -7.20%  "
Set-TextValue "D26" "64.08"
$ws.Range("E26").Value = "  -1.50%  "
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D27" "8.19"
$ws.Range("E27").Value = "  -10.62%  "
$ws.Range("B28").Value = "Bittensor"
$ws.Range("C28").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D28" "553.89"
$ws.Range("E28").Value = "  -5.48%  "
$ws.Range("D29").Value = "2.530.60"
$ws.Range("E29").Value = "  -1.25%  "
$ws.Range("D30").Value = "0.0₃0908"
$ws.Range("E30").Value = "  -2.72%  "
Set-TextValue "D31" "7.87"
$ws.Range("E31").Value = "  -0.10%  "
Set-TextValue "D32" "1.30"
$ws.Range("E32").Value = "  -5.55%  "
Set-TextValue "D33" "1.79"
$ws.Range("E33").Value = "  -4.16%  "
Set-TextValue "D34" "0.131"
$ws.Range("E34").Value = "  -1.60%  "
Set-TextValue "D35" "1.00"
$ws.Range("E35").Value = "  +0.29%  "
Set-TextValue "D36" "152.04"
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("E37").Value = "  -0.40%  "
Set-TextValue "D38" "0.365"
$ws.Range("E38").Value = "  -1.79%  "
Set-TextValue "D39" "4.49"
$ws.Range("E39").Value = "  -5.81%  "
Set-TextValue "D40" "18.04"
$ws.Range("E40").Value = "  -1.50%  "
Set-TextValue "D41" "5.05"
$ws.Range("E41").Value = "  -2.26%  "
$ws.Range("E42").Value = "  +0.02%  "
Set-TextValue "D43" "1.63"
$ws.Range("E43").Value = "  -2.61%  "
Set-TextValue "D44" "2.32"
$ws.Range("E44").Value = "  -1.92%  "
$ws.Range("D45").Value = "0.0₆0287"
$ws.Range("E45").Value = "  +2.63%  "
Set-TextValue "D46" "141.05"
$ws.Range("E46").Value = "  -0.18%  "
Set-TextValue "D47" "3.46"
$ws.Range("E47").Value = "  -2.98%  "
Set-TextValue "D48" "0.581"
$ws.Range("E48").Value = "  -2.88%  "
Set-TextValue "D49" "0.0496"
$ws.Range("E49").Value = "  -2.71%  "
Set-TextValue "D50" "18.92"
$ws.Range("E50").Value = "  -3.66%  "
Set-TextValue "D51" "0.0893"
$ws.Range("E51").Value = "  -0.57%  "
